# "nuevos experimentos no convexos" -- refresh the generated experiment
# numbers across the workbook (non-convex MorganPatrone2006a generator
# re-run with a new alpha).
#
# Every touched cell (besides Vector_Alpha!A2, a genuine number) already
# holds text -- algebraic expressions or stringified floats written by the
# generator -- so each target range is forced to Text format before the
# new value is poured in (otherwise a numeric-looking string such as
# "0.3" or "6.0" would silently be re-typed as a number), then restored to
# the "Normal" cell style so no stray formatting is left behind on the
# cells themselves.

$wb = $excel.ActiveWorkbook

# NOTE: PowerShell hashtables key-match case-insensitively by default, and
# this workbook has two sheets whose names differ only by case
# ("Vector_bf" vs "Vector_BF"), so sheet refs are kept in separate
# variables (indexed by position) rather than a hashtable keyed by name.
$wsLider    = $wb.Worksheets.Item(2)  # Restricciones_del_lider
$wsFollower = $wb.Worksheets.Item(3)  # Restricciones_del_follower
$wsPunto    = $wb.Worksheets.Item(4)  # Punto_modificado
$wsVecbf    = $wb.Worksheets.Item(5)  # Vector_bf
$wsVecBF    = $wb.Worksheets.Item(6)  # Vector_BF
$wsAlpha    = $wb.Worksheets.Item(7)  # Vector_Alpha

$textRanges = @(
    ,@($wsLider,    "A2:D3")
    ,@($wsFollower, "A2:F3")
    ,@($wsPunto,    "A2:B2")
    ,@($wsVecbf,    "A2")
    ,@($wsVecBF,    "A2:A3")
)

# Force every touched range to Text so the new strings aren't silently
# reinterpreted as numbers.
foreach ($pair in $textRanges) {
    $pair[0].Range($pair[1]).NumberFormat = "@"
}

# Restricciones_del_lider
$wsLider.Range("A2").Value = "2.8499999999999996 - x"
$wsLider.Range("B2").Value = "-3.3499999999999996"
$wsLider.Range("D2").Value = "0.3"
$wsLider.Range("A3").Value = "-2.8499999999999996 + x"
$wsLider.Range("B3").Value = "2.3499999999999996"
$wsLider.Range("D3").Value = "0.09"

# Restricciones_del_follower
$wsFollower.Range("A2").Value = "-13.271929824561402 + 2.9824561403508776y"
$wsFollower.Range("B2").Value = "12.271929824561402"
$wsFollower.Range("D2").Value = "0.19"
$wsFollower.Range("E2").Value = "6.0"
$wsFollower.Range("F2").Value = "6.800000000000001"
$wsFollower.Range("A3").Value = "-4.440892098500626e-16 + 1.1102230246251565e-16y"
$wsFollower.Range("B3").Value = "-0.9999999999999996"
$wsFollower.Range("D3").Value = "0.79"
$wsFollower.Range("E3").Value = "9.0"
$wsFollower.Range("F3").Value = "1.9"

# Punto_modificado
$wsPunto.Range("A2").Value = "2.8499999999999996"
$wsPunto.Range("B2").Value = "4.449999999999999"

# Vector_bf
$wsVecbf.Range("A2").Value = "-3.4166666666666665"

# Vector_BF
$wsVecBF.Range("A2").Value = "-1.0700000000000003"
$wsVecBF.Range("A3").Value = "-16.894736842105267"

# Vector_Alpha (genuine numeric cell, not text)
$wsAlpha.Range("A2").Value = 2.2800000000000002

# Restore the plain "Normal" style on every cell we forced to Text above,
# so only the cell contents changed -- not their formatting.
foreach ($pair in $textRanges) {
    $pair[0].Range($pair[1]).Style = "Normal"
}
